# Updates cryptos list data (Coin/Link/Price/Volume(1h)) to match the
# latest scrape, per commit "Updated cryptos list ... with GitHub Actions".
#
# Price values (column D) are written through a text formula + copy/
# paste-values round-trip for any value that Excel would otherwise auto-
# convert to a number (stripping meaningful trailing zeros / precision,
# e.g. "1.00" -> 1, "0.497" -> 0.497000000000004). This keeps the cells
# as plain text - exactly matching the source data - without touching
# any cell formatting/styles.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.813.81"
$ws.Range("E2").Value = "  -4.64%  "
$ws.Range("D3").Value = "3.211.75"
$ws.Range("E3").Value = "  -8.23%  "
$ws.Range("D4").Formula = "=""1.00"""
$ws.Range("D4").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Formula = "=""597.02"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("D6").Formula = "=""151.40"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  -12.26%  "
$ws.Range("D7").Formula = "=""0.999"""
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "3.205.64"
$ws.Range("E8").Value = "  -8.23%  "
$ws.Range("D9").Formula = "=""0.543"""
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  -10.66%  "
$ws.Range("E10").Value = "  -10.37%  "
$ws.Range("D11").Formula = "=""6.52"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  -9.30%  "
$ws.Range("D12").Formula = "=""0.497"""
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  -15.15%  "
$ws.Range("D13").Formula = "=""39.02"""
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  -15.25%  "
$ws.Range("E14").Value = "  -11.31%  "
$ws.Range("D15").Value = "3.732.04"
$ws.Range("E15").Value = "  -8.30%  "
$ws.Range("D16").Value = "66.839.54"
$ws.Range("E16").Value = "  -4.64%  "
$ws.Range("D17").Value = "3.216.27"
$ws.Range("E17").Value = "  -8.26%  "
$ws.Range("E18").Value = "  -4.36%  "
$ws.Range("D19").Formula = "=""533.33"""
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  -12.85%  "
$ws.Range("D20").Formula = "=""7.15"""
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  -14.33%  "
$ws.Range("D21").Formula = "=""14.96"""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  -14.65%  "
$ws.Range("D22").Formula = "=""0.762"""
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  -13.08%  "
$ws.Range("E23").Value = "  -12.71%  "
$ws.Range("D24").Formula = "=""13.85"""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -11.09%  "
$ws.Range("D25").Formula = "=""85.42"""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").Formula = "=""3.19"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  -13.99%  "
$ws.Range("E28").Value = "  -14.10%  "
$ws.Range("E29").Value = "  -8.85%  "
$ws.Range("D30").Formula = "=""29.20"""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  -14.04%  "
$ws.Range("D31").Formula = "=""2.65"""
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  -10.47%  "
$ws.Range("E32").Value = "  -10.17%  "
$ws.Range("D33").Formula = "=""549.89"""
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  -12.55%  "
$ws.Range("D34").Formula = "=""6.56"""
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  -18.24%  "
$ws.Range("E35").Value = "  -15.81%  "
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("D37").Formula = "=""53.47"""
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  -5.70%  "
$ws.Range("D38").Formula = "=""0.0431"""
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  -9.43%  "
$ws.Range("D39").Formula = "=""0.0867"""
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  -12.69%  "
$ws.Range("E40").Value = "  -12.39%  "
$ws.Range("E41").Value = "  -12.03%  "
$ws.Range("D42").Value = "2.914.38"
$ws.Range("E42").Value = "  -13.20%  "
$ws.Range("D43").Formula = "=""2.67"""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -23.37%  "
$ws.Range("D44").Formula = "=""0.264"""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  -14.61%  "
$ws.Range("D45").Value = "0.0₃0585"
$ws.Range("E45").Value = "  -19.82%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Formula = "=""26.72"""
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  -16.08%  "
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").Formula = "=""2.40"""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  -17.15%  "
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").Formula = "=""2.13"""
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  -16.06%  "
$ws.Range("E50").Value = "  -12.24%  "
$ws.Range("D51").Formula = "=""121.39"""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  -8.66%  "

# Clear the clipboard marching-ants state left over from the paste-values steps.
$excel.CutCopyMode = 0
